$wb = $excel.ActiveWorkbook

# Mapping of row -> [old value, new value] for column F ("想去人数")
$changes = @{
    3  = 143
    4  = 1388
    5  = 1622
    9  = 196
    11 = 76
    13 = 291
    14 = 333
    16 = 1798
    20 = 708
    23 = 4335
    25 = 307
    26 = 1153
    27 = 504
    29 = 699
    31 = 351
    33 = 181
}

# Both "展览" and "全部类型" sheets contain the same data and need updating.
$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $changes.Keys) {
        $ws.Cells.Item($row, 6).Value = $changes[$row]
    }
}
